$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- 1. Update the "Status" text: the handback is no longer in sync with en-US ---
$newStatus = "Handed back: not in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- 2. New handback report was generated for row 3 (6f2007ba... file) on both
#        locale sheets, so its "Correspond Handback DateTime" gets a fresh value ---
$wsZhCn.Range("K3").Value = "2016-10-18 13:05:27"
$wsDeDe.Range("K3").Value = "2016-10-18 13:05:45"

# --- 3. The Status column grew wider to fit the new, longer status text ---
$wsOverview.Range("E1").ColumnWidth = 32.6
$wsOverview.Range("F1").ColumnWidth = 32.6
$wsZhCn.Range("C1").ColumnWidth = 32.6
$wsDeDe.Range("C1").ColumnWidth = 32.6
